$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new quotation row (row 19) with the latest fund prices.
$ws.Range("A19").Value = 45923
$ws.Range("B19").Value = "20,8632"
$ws.Range("C19").Value = "15,1509"
$ws.Range("D19").Value = "14,9245"
$ws.Range("E19").Value = "14,9245"

# Match the date formatting/style used by the rest of column A.
$ws.Range("A19").NumberFormat = $ws.Range("A18").NumberFormat
